$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows 9-22 which no longer exist in the new layout
$ws.Range("A9:A22").ClearContents()

# Update the remaining rows with the combined tuple-style text
$ws.Range("A2").Value = "('The Avenger', ['Hero', '{3}, {T}: Target creature you control gains deathtouch until end of turn.'])"
$ws.Range("A3").Value = "('The Harvester', ['Hero', '{T}: Draw a card, then discard a card.'])"
$ws.Range("A4").Value = "('The Hunter', ['Hero', '{T}: Target creature you control gets +1/+1 until end of turn.'])"
$ws.Range("A5").Value = "('The Philosopher', ['Hero', '{2}, {T}: Tap target creature.'])"
$ws.Range("A6").Value = "('The Protector', ['Hero', '{T}: Prevent the next 1 damage that would be dealt to any target this turn.'])"
$ws.Range("A7").Value = "('The Slayer', ['Hero', 'You start the game with an additional 7 life.'])"
$ws.Range("A8").Value = "('The Warrior', ['Hero', '{T}: Target creature you control gains haste until end of turn.'])"
